$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.91
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 1.42
$ws.Range("K2").Value = 2.47
$ws.Range("L2").Value = 2.18
$ws.Range("M2").Value = 1.53
$ws.Range("N2").Value = 1.5
$ws.Range("O2").Value = 2.27
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 5.8
$ws.Range("S2").Value = 8
$ws.Range("T2").Value = 8.75
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 17.5
$ws.Range("X2").Value = 7.4
$ws.Range("Y2").Value = 6.4
$ws.Range("Z2").Value = 18.5
$ws.Range("AA2").Value = 110
$ws.Range("AB2").Value = 9
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 14
$ws.Range("AE2").Value = 60
$ws.Range("AF2").Value = 45
$ws.Range("AG2").Value = 60
$ws.Range("G4").Value = 2.2
$ws.Range("I4").Value = 3.2
$ws.Range("T4").Value = 9
$ws.Range("X4").Value = 9.5
$ws.Range("AD4").Value = 12
$ws.Range("G6").Value = 2.55
$ws.Range("I6").Value = 2.6
$ws.Range("L6").Value = 1.91
$ws.Range("M6").Value = 1.91
$ws.Range("AB6").Value = 9
$ws.Range("J8").Value = 1.5
$ws.Range("K8").Value = 2.63
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 3.1
$ws.Range("S10").Value = 11
$ws.Range("AI10").Value = 1.07
$ws.Range("AJ10").Value = 9
$ws.Range("J11").Value = 1.25
$ws.Range("K11").Value = 4
$ws.Range("P11").Value = 1.73
$ws.Range("Q11").Value = 2
$ws.Range("Y11").Value = 7
$ws.Range("AH11").Value = 201
$ws.Range("AI11").Value = 1.05
$ws.Range("AJ11").Value = 11
$ws.Range("G13").Value = 2.7
$ws.Range("I13").Value = 2.6
$ws.Range("L13").Value = 2.4
$ws.Range("M13").Value = 1.53
$ws.Range("P13").Value = 2
$ws.Range("Q13").Value = 1.75
$ws.Range("X13").Value = 7.5
$ws.Range("Y13").Value = 6
$ws.Range("Z13").Value = 17
$ws.Range("AB13").Value = 7
$ws.Range("AC13").Value = 12
$ws.Range("AD13").Value = 11
$ws.Range("AE13").Value = 26
$ws.Range("AH13").Value = 451
$ws.Range("G14").Value = 2.7
$ws.Range("I14").Value = 2.6
$ws.Range("R14").Value = 7.5
$ws.Range("T14").Value = 11
$ws.Range("U14").Value = 29
$ws.Range("Y14").Value = 6
$ws.Range("AA14").Value = 51
$ws.Range("AB14").Value = 7.5
$ws.Range("AH14").Value = 351
$ws.Range("AI14").Value = 1.08
$ws.Range("AJ14").Value = 8
$ws.Range("AJ15").Value = 10
$ws.Range("AI17").Value = 1.06
$ws.Range("AJ17").Value = 10
$ws.Range("AB19").Value = 18.5
$ws.Range("H20").Value = 6.6
$ws.Range("I20").Value = 18.5
$ws.Range("P20").Value = 2.1
$ws.Range("Q20").Value = 1.66
$ws.Range("R20").Value = 8.5
$ws.Range("T20").Value = 9.25
$ws.Range("X20").Value = 18.5
$ws.Range("Z20").Value = 28
$ws.Range("AA20").Value = 110
$ws.Range("AB20").Value = 50
$ws.Range("AC20").Value = 175
$ws.Range("AD20").Value = 55
$ws.Range("AF20").Value = 300
$ws.Range("AG20").Value = 150
$ws.Range("G21").Value = 6.1
$ws.Range("H21").Value = 4.1
$ws.Range("I21").Value = 1.42
$ws.Range("L21").Value = 1.53
$ws.Range("M21").Value = 2.18
$ws.Range("P21").Value = 1.71
$ws.Range("Q21").Value = 2.02
$ws.Range("R21").Value = 17
$ws.Range("S21").Value = 35
$ws.Range("T21").Value = 15.5
$ws.Range("U21").Value = 100
$ws.Range("V21").Value = 45
$ws.Range("W21").Value = 37
$ws.Range("X21").Value = 13.5
$ws.Range("Y21").Value = 7.4
$ws.Range("Z21").Value = 12.5
$ws.Range("AB21").Value = 7
$ws.Range("AC21").Value = 6.5
$ws.Range("AD21").Value = 6.9
$ws.Range("AE21").Value = 8.75
$ws.Range("AF21").Value = 9.25
$ws.Range("AH21").Value = 250
$ws.Range("S22").Value = 8
$ws.Range("AF22").Value = 41
$ws.Range("AI22").Value = 1.05
$ws.Range("AJ22").Value = 11
$ws.Range("G23").Value = 3.4
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 2
$ws.Range("L23").Value = 1.85
$ws.Range("P23").Value = 1.69
$ws.Range("S23").Value = 19
$ws.Range("U23").Value = 41
$ws.Range("V23").Value = 29
$ws.Range("Z23").Value = 15
$ws.Range("AC23").Value = 10
$ws.Range("P24").Value = 1.8
$ws.Range("Q24").Value = 1.8
$ws.Range("AE24").Value = 51
$ws.Range("L25").Value = 1.85
$ws.Range("M25").Value = 1.95
$ws.Range("P25").Value = 1.69
$ws.Range("M26").Value = 1.95
$ws.Range("N26").Value = 1.36
$ws.Range("O26").Value = 3
$ws.Range("P26").Value = 1.63
$ws.Range("R26").Value = 9.5
$ws.Range("W26").Value = 26
$ws.Range("X26").Value = 12
$ws.Range("AG26").Value = 26
$ws.Range("AH26").Value = 151
$ws.Range("H27").Value = 7.1
$ws.Range("I27").Value = 24
$ws.Range("W27").Value = 28
$ws.Range("AB27").Value = 75
$ws.Range("AG27").Value = 175
$ws.Range("G28").Value = 1.41
$ws.Range("H28").Value = 4.7
$ws.Range("I28").Value = 6.2
$ws.Range("J28").Value = 1.16
$ws.Range("K28").Value = 4.9
$ws.Range("L28").Value = 1.51
$ws.Range("M28").Value = 2.46
$ws.Range("P28").Value = 1.7
$ws.Range("Q28").Value = 2.07
$ws.Range("R28").Value = 7.2
$ws.Range("S28").Value = 6.4
$ws.Range("T28").Value = 6.6
$ws.Range("U28").Value = 8
$ws.Range("V28").Value = 8.4
$ws.Range("W28").Value = 17
$ws.Range("X28").Value = 15
$ws.Range("Y28").Value = 7.2
$ws.Range("Z28").Value = 13
$ws.Range("AA28").Value = 50
$ws.Range("AB28").Value = 17
$ws.Range("AC28").Value = 35
$ws.Range("AD28").Value = 15
$ws.Range("AE28").Value = 101
$ws.Range("AF28").Value = 50
$ws.Range("AG28").Value = 40
$ws.Range("AH28").Value = 101
$ws.Range("G29").Value = 1.38
$ws.Range("H29").Value = 4.2
$ws.Range("I29").Value = 8.6
$ws.Range("J29").Value = 1.36
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = 2.09
$ws.Range("M29").Value = 1.69
$ws.Range("N29").Value = 1.4
$ws.Range("O29").Value = 2.75
$ws.Range("P29").Value = 2.49
$ws.Range("Q29").Value = 1.5
$ws.Range("R29").Value = 4.1
$ws.Range("S29").Value = 4.3
$ws.Range("T29").Value = 7.2
$ws.Range("U29").Value = 6.4
$ws.Range("V29").Value = 11
$ws.Range("W29").Value = 40
$ws.Range("X29").Value = 7
$ws.Range("Y29").Value = 6.8
$ws.Range("Z29").Value = 24
$ws.Range("AA29").Value = 101
$ws.Range("AB29").Value = 13
$ws.Range("AC29").Value = 45
$ws.Range("AD29").Value = 24
$ws.Range("AE29").Value = 101
$ws.Range("AF29").Value = 101
$ws.Range("AG29").Value = 101
$ws.Range("AH29").Value = 101
$ws.Range("S30").Value = 10
$ws.Range("AD30").Value = 12
$ws.Range("AH30").Value = 301
$ws.Range("AI30").Value = 1.07
$ws.Range("AJ30").Value = 9
$ws.Range("J31").Value = 1.29
$ws.Range("K31").Value = 3.5
$ws.Range("L31").Value = 1.93
$ws.Range("M31").Value = 1.93
$ws.Range("N31").Value = 1.36
$ws.Range("O31").Value = 3
$ws.Range("V31").Value = 13
$ws.Range("X31").Value = 11
$ws.Range("AH31").Value = 251
$ws.Range("I32").Value = 6.5
$ws.Range("J32").Value = 1.3
$ws.Range("K32").Value = 3.4
$ws.Range("L32").Value = 2
$ws.Range("M32").Value = 1.85
$ws.Range("P32").Value = 2.2
$ws.Range("Q32").Value = 1.62
$ws.Range("R32").Value = 6
$ws.Range("W32").Value = 34
$ws.Range("X32").Value = 9.5
$ws.Range("AJ32").Value = 9.5
$ws.Range("H33").Value = 2.95
$ws.Range("I33").Value = 2.67
$ws.Range("J33").Value = 1.38
$ws.Range("K33").Value = 2.82
$ws.Range("O33").Value = 2.65
$ws.Range("R33").Value = 8.25
$ws.Range("X33").Value = 6.3
$ws.Range("Y33").Value = 5.7
$ws.Range("Z33").Value = 13.5
$ws.Range("AD33").Value = 10
$ws.Range("AE33").Value = 32
$ws.Range("AF33").Value = 25
$ws.Range("AI33").Value = 1.09
$ws.Range("AJ33").Value = 6.3
